# Apply updated leve-profit figures (currentAveragePrice / LevePrice / LeveProfit
# columns H:N) across the per-job tables, as refreshed by the scheduled pricing
# runner. Cells that did not exist before (newly computed LeveProfit) are created;
# cells that no longer apply (e.g. HQ profit column when HQ price data dropped)
# are cleared so they stay absent, matching the source feed output.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 274.70587
$ws.Range("I2").Value = 213.07692
$ws.Range("J2").Value = 475
$ws.Range("K2").Value = 213.07692
$ws.Range("L2").Value = 475
$ws.Range("M2").Value = -100.07692
$ws.Range("N2").Value = -701
$ws.Range("H33").Value = 165.6
$ws.Range("I33").Value = 172
$ws.Range("J33").Value = 140
$ws.Range("K33").Value = 172
$ws.Range("L33").Value = 140
$ws.Range("M33").Value = 57
$ws.Range("N33").Value = -598
$ws.Range("H75").Value = 38580
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 38580
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 38580
$ws.Range("N75").Value = -40452
$ws.Range("H78").Value = 38580
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 38580
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 115740
$ws.Range("N78").Value = -125100
$ws.Range("H113").Value = 2730.6924
$ws.Range("I113").Value = 2775
$ws.Range("J113").Value = 2711
$ws.Range("K113").Value = 2775
$ws.Range("L113").Value = 2711
$ws.Range("M113").Value = 479
$ws.Range("N113").Value = -9219
$ws.Range("H137").Value = 3477.25
$ws.Range("I137").Value = 3134.5789
$ws.Range("J137").Value = 4200.6665
$ws.Range("K137").Value = 9403.736699999999
$ws.Range("L137").Value = 12601.9995
$ws.Range("M137").Value = -6853.736699999999
$ws.Range("N137").Value = -17701.9995
$ws.Range("H138").Value = 2006.2162
$ws.Range("I138").Value = 1695.7241
$ws.Range("J138").Value = 3131.75
$ws.Range("K138").Value = 5087.1723
$ws.Range("L138").Value = 9395.25
$ws.Range("M138").Value = 52.82769999999982
$ws.Range("N138").Value = -19675.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 583033.25
$ws.Range("I32").Value = 631083.25
$ws.Range("J32").Value = 35262.8
$ws.Range("K32").Value = 631083.25
$ws.Range("L32").Value = 35262.8
$ws.Range("M32").Value = -630796.25
$ws.Range("N32").Value = -35836.8
$ws.Range("H61").Value = 1814.7142
$ws.Range("I61").Value = 1504.766
$ws.Range("J61").Value = 3433.3333
$ws.Range("K61").Value = 1504.766
$ws.Range("L61").Value = 3433.3333
$ws.Range("M61").Value = -1292.766
$ws.Range("N61").Value = -3857.3333
$ws.Range("H74").Value = 1339.0385
$ws.Range("I74").Value = 904.5789
$ws.Range("J74").Value = 2518.2856
$ws.Range("K74").Value = 904.5789
$ws.Range("L74").Value = 2518.2856
$ws.Range("M74").Value = -30.57889999999998
$ws.Range("N74").Value = -4266.2856
$ws.Range("H77").Value = 1339.0385
$ws.Range("I77").Value = 904.5789
$ws.Range("J77").Value = 2518.2856
$ws.Range("K77").Value = 4522.8945
$ws.Range("L77").Value = 12591.428
$ws.Range("M77").Value = -154.8945000000003
$ws.Range("N77").Value = -21327.428
$ws.Range("H86").Value = 111137780
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 111137780
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 111137780
$ws.Range("N86").Value = -111140152
$ws.Range("H89").Value = 111137780
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 111137780
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 333413340
$ws.Range("N89").Value = -333425196
$ws.Range("H101").Value = 79602
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 79602
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 79602
$ws.Range("N101").Value = -86092
$ws.Range("H132").Value = 4130.1396
$ws.Range("I132").Value = 2745.75
$ws.Range("J132").Value = 5878.8423
$ws.Range("K132").Value = 8237.25
$ws.Range("L132").Value = 17636.5269
$ws.Range("M132").Value = -5707.25
$ws.Range("N132").Value = -22696.5269
$ws.Range("H136").Value = 1814.7142
$ws.Range("I136").Value = 1504.766
$ws.Range("J136").Value = 3433.3333
$ws.Range("K136").Value = 4514.298000000001
$ws.Range("L136").Value = 10299.9999
$ws.Range("M136").Value = -1964.298000000001
$ws.Range("N136").Value = -15399.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2075
$ws.Range("I94").Value = 1300
$ws.Range("J94").Value = 2333.3333
$ws.Range("K94").Value = 1300
$ws.Range("L94").Value = 2333.3333
$ws.Range("M94").Value = -849
$ws.Range("N94").Value = -3235.3333
$ws.Range("H105").Value = 3272
$ws.Range("I105").Value = 3272
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3272
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -1525
$ws.Range("H134").Value = 2387.739
$ws.Range("I134").Value = 1844.5454
$ws.Range("J134").Value = 2885.6667
$ws.Range("K134").Value = 5533.6362
$ws.Range("L134").Value = 8657.000100000001
$ws.Range("M134").Value = -2998.6362
$ws.Range("N134").Value = -13727.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2285.5715
$ws.Range("I105").Value = 2000
$ws.Range("J105").Value = 2666.3333
$ws.Range("K105").Value = 2000
$ws.Range("L105").Value = 2666.3333
$ws.Range("M105").Value = -253
$ws.Range("N105").Value = -6160.3333
$ws.Range("H112").Value = 36666.5
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 36666.5
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 36666.5
$ws.Range("N112").Value = -39620.5
$ws.Range("H118").Value = 44444
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 44444
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 44444
$ws.Range("N118").Value = -47758
$ws.Range("H132").Value = 1750.9375
$ws.Range("I132").Value = 1354.1904
$ws.Range("J132").Value = 2508.3635
$ws.Range("K132").Value = 4062.5712
$ws.Range("L132").Value = 7525.0905
$ws.Range("M132").Value = -1532.5712
$ws.Range("N132").Value = -12585.0905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 780.2222
$ws.Range("I5").Value = 543.7692
$ws.Range("J5").Value = 1395
$ws.Range("K5").Value = 1631.3076
$ws.Range("L5").Value = 4185
$ws.Range("M5").Value = -1519.3076
$ws.Range("N5").Value = -4409
$ws.Range("H80").Value = 4185.8887
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 4185.8887
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 12557.6661
$ws.Range("N80").Value = -14429.6661
$ws.Range("H83").Value = 4185.8887
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 4185.8887
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 37672.99830000001
$ws.Range("N83").Value = -47032.99830000001
$ws.Range("H113").Value = 799.55316
$ws.Range("I113").Value = 570.4545000000001
$ws.Range("J113").Value = 1001.16
$ws.Range("K113").Value = 1711.3635
$ws.Range("L113").Value = 3003.48
$ws.Range("M113").Value = 458.6364999999998
$ws.Range("N113").Value = -7343.48
$ws.Range("H135").Value = 780.2222
$ws.Range("I135").Value = 543.7692
$ws.Range("J135").Value = 1395
$ws.Range("K135").Value = 4893.922799999999
$ws.Range("L135").Value = 12555
$ws.Range("M135").Value = -2358.922799999999
$ws.Range("N135").Value = -17625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 56601770
$ws.Range("I80").Value = 84835000
$ws.Range("J80").Value = 135296
$ws.Range("K80").Value = 84835000
$ws.Range("L80").Value = 135296
$ws.Range("M80").Value = -84834002
$ws.Range("N80").Value = -137292
$ws.Range("H83").Value = 56601770
$ws.Range("I83").Value = 84835000
$ws.Range("J83").Value = 135296
$ws.Range("K83").Value = 424175000
$ws.Range("L83").Value = 676480
$ws.Range("M83").Value = -424170008
$ws.Range("N83").Value = -686464
$ws.Range("H102").Value = 1201
$ws.Range("I102").Value = 1238.4
$ws.Range("J102").Value = 1014
$ws.Range("K102").Value = 1238.4
$ws.Range("L102").Value = 1014
$ws.Range("M102").Value = 383.5999999999999
$ws.Range("N102").Value = -4258
$ws.Range("H122").Value = 2087.2942
$ws.Range("I122").Value = 1955.875
$ws.Range("J122").Value = 2402.7
$ws.Range("K122").Value = 5867.625
$ws.Range("L122").Value = 7208.099999999999
$ws.Range("M122").Value = -3417.625
$ws.Range("N122").Value = -12108.1
$ws.Range("H126").Value = 4497.5
$ws.Range("I126").Value = 4497.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 13492.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -11022.5
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 3659.9048
$ws.Range("I132").Value = 2842.1538
$ws.Range("J132").Value = 4988.75
$ws.Range("K132").Value = 8526.4614
$ws.Range("L132").Value = 14966.25
$ws.Range("M132").Value = -5996.4614
$ws.Range("N132").Value = -20026.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 14287484
$ws.Range("I16").Value = 1633.6666
$ws.Range("J16").Value = 35716260
$ws.Range("K16").Value = 1633.6666
$ws.Range("L16").Value = 35716260
$ws.Range("M16").Value = -1463.6666
$ws.Range("N16").Value = -35716600
$ws.Range("H55").Value = 919.1539
$ws.Range("I55").Value = 821.125
$ws.Range("J55").Value = 1076
$ws.Range("K55").Value = 821.125
$ws.Range("L55").Value = 1076
$ws.Range("M55").Value = -648.125
$ws.Range("N55").Value = -1422
$ws.Range("H122").Value = 3328.125
$ws.Range("I122").Value = 2490
$ws.Range("J122").Value = 3709.0908
$ws.Range("K122").Value = 7470
$ws.Range("L122").Value = 11127.2724
$ws.Range("M122").Value = -5020
$ws.Range("N122").Value = -16027.2724
$ws.Range("H132").Value = 2097.2666
$ws.Range("I132").Value = 1432.4762
$ws.Range("J132").Value = 3648.4443
$ws.Range("K132").Value = 4297.4286
$ws.Range("L132").Value = 10945.3329
$ws.Range("M132").Value = -1767.4286
$ws.Range("N132").Value = -16005.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
